# Auto-generated cell value updates for Midgardsormr_Profits sheets
# (FFXIV crafting-leve profit data refresh from scheduled runner)
$wb = $excel.ActiveWorkbook

# --- ALC sheet ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H6").Value = 796.73334
$ws.Range("I6").Value = 401.6
$ws.Range("K6").Value = 1204.8
$ws.Range("M6").Value = -1092.8
$ws.Range("H29").Value = 6670.8
$ws.Range("J29").Value = 10501.333
$ws.Range("L29").Value = 31503.999
$ws.Range("N29").Value = -32065.999
$ws.Range("H38").Value = 1772.6
$ws.Range("I38").Value = 963.5
$ws.Range("K38").Value = 2890.5
$ws.Range("M38").Value = -2518.5
$ws.Range("H39").Value = 412.13333
$ws.Range("I39").Value = 321.69232
$ws.Range("J39").Value = 1000
$ws.Range("K39").Value = 965.07696
$ws.Range("L39").Value = 3000
$ws.Range("M39").Value = -669.07696
$ws.Range("N39").Value = -3592
$ws.Range("H62").Value = 7100.3335
$ws.Range("I62").Value = 7284.5835
$ws.Range("J62").Value = 6916.0835
$ws.Range("K62").Value = 7284.5835
$ws.Range("L62").Value = 6916.0835
$ws.Range("M62").Value = -6660.5835
$ws.Range("N62").Value = -8164.0835
$ws.Range("H65").Value = 7100.3335
$ws.Range("I65").Value = 7284.5835
$ws.Range("J65").Value = 6916.0835
$ws.Range("K65").Value = 36422.9175
$ws.Range("L65").Value = 34580.4175
$ws.Range("M65").Value = -33302.9175
$ws.Range("N65").Value = -40820.4175
$ws.Range("H74").Value = 3106.5
$ws.Range("I74").Value = 3106.5
$ws.Range("K74").Value = 3106.5
$ws.Range("M74").Value = -2170.5
$ws.Range("H77").Value = 3106.5
$ws.Range("I77").Value = 3106.5
$ws.Range("K77").Value = 15532.5
$ws.Range("M77").Value = -10852.5
$ws.Range("H86").Value = 2962.2778
$ws.Range("I86").Value = 3088.625
$ws.Range("J86").Value = 2861.2
$ws.Range("K86").Value = 3088.625
$ws.Range("L86").Value = 2861.2
$ws.Range("M86").Value = -1965.625
$ws.Range("N86").Value = -5107.2
$ws.Range("H89").Value = 2962.2778
$ws.Range("I89").Value = 3088.625
$ws.Range("J89").Value = 2861.2
$ws.Range("K89").Value = 15443.125
$ws.Range("L89").Value = 14306
$ws.Range("M89").Value = -9827.125
$ws.Range("N89").Value = -25538
$ws.Range("H106").Value = 2230.12
$ws.Range("I106").Value = 1788.6
$ws.Range("J106").Value = 3996.2
$ws.Range("K106").Value = 1788.6
$ws.Range("L106").Value = 3996.2
$ws.Range("M106").Value = -1157.6
$ws.Range("N106").Value = -5258.2
$ws.Range("H137").Value = 8665.02
$ws.Range("I137").Value = 19024.611
$ws.Range("J137").Value = 3014.3333
$ws.Range("K137").Value = 57073.833
$ws.Range("L137").Value = 9042.999899999999
$ws.Range("M137").Value = -54523.833
$ws.Range("N137").Value = -14142.9999
$ws.Range("H138").Value = 1702.6666
$ws.Range("I138").Value = 1702.6666
$ws.Range("K138").Value = 5107.9998
$ws.Range("M138").Value = 32.0002000000004

# --- ARM sheet ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1806.697
$ws.Range("I2").Value = 1412.9445
$ws.Range("K2").Value = 1412.9445
$ws.Range("M2").Value = -1299.9445
$ws.Range("H61").Value = 1399.88
$ws.Range("I61").Value = 1086.8695
$ws.Range("J61").Value = 4999.5
$ws.Range("K61").Value = 1086.8695
$ws.Range("L61").Value = 4999.5
$ws.Range("M61").Value = -874.8695
$ws.Range("N61").Value = -5423.5
$ws.Range("H74").Value = 858560.1
$ws.Range("I74").Value = 858560.1
$ws.Range("J74").Value = 0
$ws.Range("K74").Value = 858560.1
$ws.Range("L74").Value = 0
$ws.Range("M74").Value = -857686.1
$ws.Range("H77").Value = 858560.1
$ws.Range("I77").Value = 858560.1
$ws.Range("J77").Value = 0
$ws.Range("K77").Value = 4292800.5
$ws.Range("L77").Value = 0
$ws.Range("M77").Value = -4288432.5
$ws.Range("H97").Value = 1626.48
$ws.Range("I97").Value = 1243.7812
$ws.Range("J97").Value = 2306.8333
$ws.Range("K97").Value = 1243.7812
$ws.Range("L97").Value = 2306.8333
$ws.Range("M97").Value = -747.7811999999999
$ws.Range("N97").Value = -3298.8333
$ws.Range("H102").Value = 5253.032
$ws.Range("I102").Value = 4806.1304
$ws.Range("K102").Value = 4806.1304
$ws.Range("M102").Value = -3184.1304
$ws.Range("H116").Value = 1806.697
$ws.Range("I116").Value = 1412.9445
$ws.Range("K116").Value = 1412.9445
$ws.Range("M116").Value = 881.0554999999999
$ws.Range("H132").Value = 1891.279
$ws.Range("I132").Value = 1339.3529
$ws.Range("J132").Value = 3976.3333
$ws.Range("K132").Value = 4018.0587
$ws.Range("L132").Value = 11928.9999
$ws.Range("M132").Value = -1488.0587
$ws.Range("N132").Value = -16988.9999
$ws.Range("H136").Value = 1399.88
$ws.Range("I136").Value = 1086.8695
$ws.Range("J136").Value = 4999.5
$ws.Range("K136").Value = 3260.6085
$ws.Range("L136").Value = 14998.5
$ws.Range("M136").Value = -710.6085000000003
$ws.Range("N136").Value = -20098.5
$ws.Range("N74").ClearContents()
$ws.Range("N77").ClearContents()

# --- BSM sheet ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1806.697
$ws.Range("I3").Value = 1412.9445
$ws.Range("K3").Value = 1412.9445
$ws.Range("M3").Value = -1298.9445
$ws.Range("H134").Value = 3905.25
$ws.Range("I134").Value = 3778.3333
$ws.Range("J134").Value = 5047.5
$ws.Range("K134").Value = 11334.9999
$ws.Range("L134").Value = 15142.5
$ws.Range("M134").Value = -8799.999899999999
$ws.Range("N134").Value = -20212.5

# --- CRP sheet ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2633813.5
$ws.Range("I31").Value = 3032552.8
$ws.Range("J31").Value = 2135.2
$ws.Range("K31").Value = 3032552.8
$ws.Range("L31").Value = 2135.2
$ws.Range("M31").Value = -3032257.8
$ws.Range("N31").Value = -2725.2
$ws.Range("H34").Value = 2633813.5
$ws.Range("I34").Value = 3032552.8
$ws.Range("J34").Value = 2135.2
$ws.Range("K34").Value = 3032552.8
$ws.Range("L34").Value = 2135.2
$ws.Range("M34").Value = -3032350.8
$ws.Range("N34").Value = -2539.2
$ws.Range("H132").Value = 17881.777
$ws.Range("I132").Value = 19829.594
$ws.Range("J132").Value = 2299.25
$ws.Range("K132").Value = 59488.78200000001
$ws.Range("L132").Value = 6897.75
$ws.Range("M132").Value = -56958.78200000001
$ws.Range("N132").Value = -11957.75
$ws.Range("H134").Value = 4136.615
$ws.Range("I134").Value = 2915.8333
$ws.Range("K134").Value = 8747.499899999999
$ws.Range("M134").Value = -6212.499899999999

# --- CUL sheet ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 42722704
$ws.Range("I4").Value = 28640186
$ws.Range("J4").Value = 83092590
$ws.Range("K4").Value = 85920558
$ws.Range("L4").Value = 249277770
$ws.Range("M4").Value = -85920446
$ws.Range("N4").Value = -249277994
$ws.Range("H12").Value = 217.33333
$ws.Range("J12").Value = 192.25
$ws.Range("L12").Value = 576.75
$ws.Range("N12").Value = -922.75
$ws.Range("H64").Value = 0
$ws.Range("J64").Value = 0
$ws.Range("L64").Value = 0
$ws.Range("H67").Value = 0
$ws.Range("J67").Value = 0
$ws.Range("L67").Value = 0
$ws.Range("H68").Value = 4474.9614
$ws.Range("J68").Value = 4810.375
$ws.Range("L68").Value = 14431.125
$ws.Range("N68").Value = -16053.125
$ws.Range("H71").Value = 4474.9614
$ws.Range("J71").Value = 4810.375
$ws.Range("L71").Value = 43293.375
$ws.Range("N71").Value = -51405.375
$ws.Range("H107").Value = 1387.7894
$ws.Range("I107").Value = 1723.8889
$ws.Range("K107").Value = 5171.6667
$ws.Range("M107").Value = -3251.6667
$ws.Range("H131").Value = 123244.17
$ws.Range("I131").Value = 213546.84
$ws.Range("J131").Value = 2840.6
$ws.Range("K131").Value = 640640.52
$ws.Range("L131").Value = 8521.799999999999
$ws.Range("M131").Value = -635600.52
$ws.Range("N131").Value = -18601.8
$ws.Range("H137").Value = 2929.0908
$ws.Range("J137").Value = 2613.4443
$ws.Range("L137").Value = 7840.3329
$ws.Range("N137").Value = -18040.3329
$ws.Range("N64").ClearContents()
$ws.Range("N67").ClearContents()

# --- GSM sheet ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H107").Value = 277.9
$ws.Range("J107").Value = 463.75
$ws.Range("L107").Value = 463.75
$ws.Range("N107").Value = -4303.75
$ws.Range("H131").Value = 50000
$ws.Range("J131").Value = 50000
$ws.Range("L131").Value = 50000
$ws.Range("N131").Value = -60080
$ws.Range("H132").Value = 1330.2727
$ws.Range("I132").Value = 963.3
$ws.Range("J132").Value = 5000
$ws.Range("K132").Value = 2889.9
$ws.Range("L132").Value = 15000
$ws.Range("M132").Value = -359.8999999999996
$ws.Range("N132").Value = -20060
$ws.Range("H133").Value = 66430.71000000001
$ws.Range("J133").Value = 66430.71000000001
$ws.Range("L133").Value = 66430.71000000001
$ws.Range("N133").Value = -76550.71000000001

# --- LTW sheet ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 3321.889
$ws.Range("I7").Value = 3399.7144
$ws.Range("K7").Value = 3399.7144
$ws.Range("M7").Value = -3287.7144
$ws.Range("H126").Value = 3321.889
$ws.Range("I126").Value = 3399.7144
$ws.Range("K126").Value = 10199.1432
$ws.Range("M126").Value = -7729.143199999999
$ws.Range("H132").Value = 4586.6875
$ws.Range("I132").Value = 3849.375
$ws.Range("J132").Value = 5324
$ws.Range("K132").Value = 11548.125
$ws.Range("L132").Value = 15972
$ws.Range("M132").Value = -9018.125
$ws.Range("N132").Value = -21032

# --- WVR sheet ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H64").Value = 92247.5
$ws.Range("J64").Value = 92333.336
$ws.Range("L64").Value = 92333.336
$ws.Range("N64").Value = -92829.336
$ws.Range("H67").Value = 92247.5
$ws.Range("J67").Value = 92333.336
$ws.Range("L67").Value = 92333.336
$ws.Range("N67").Value = -94049.336
$ws.Range("H122").Value = 81481.17
$ws.Range("I122").Value = 93430.16
$ws.Range("K122").Value = 280290.48
$ws.Range("M122").Value = -277840.48
$ws.Range("H132").Value = 23580.8
$ws.Range("I132").Value = 26976.166
$ws.Range("K132").Value = 80928.49800000001
$ws.Range("M132").Value = -78398.49800000001
$ws.Range("H138").Value = 99999.5
$ws.Range("J138").Value = 99999.5
$ws.Range("L138").Value = 99999.5
$ws.Range("N138").Value = -110279.5
